# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to the latest scraped values.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new value
$exhibitUpdates = @{
    "F2"  = 4
    "F3"  = 12695
    "F4"  = 24
    "F6"  = 48
    "F8"  = 13
    "F9"  = 2
    "F10" = 12592
    "F11" = 254
    "F12" = 4
    "F13" = 4935
    "F14" = 4982
    "F15" = 169
    "F16" = 79
    "F17" = 429
    "F22" = 371
    "F24" = 80
    "F25" = 5219
}

foreach ($addr in $exhibitUpdates.Keys) {
    $wsExhibit.Range($addr).Value = $exhibitUpdates[$addr]
}

# Sheet "全部类型": row -> new value
$allUpdates = @{
    "F3"  = 4
    "F4"  = 12695
    "F5"  = 24
    "F7"  = 48
    "F9"  = 13
    "F10" = 2
    "F11" = 12592
    "F12" = 254
    "F13" = 4
    "F14" = 4935
    "F15" = 4983
    "F16" = 169
    "F17" = 79
    "F18" = 429
    "F23" = 371
    "F25" = 80
    "F26" = 5219
}

foreach ($addr in $allUpdates.Keys) {
    $wsAll.Range($addr).Value = $allUpdates[$addr]
}

$wb.Save()
